# Flexible tables in PDF and configurations to include/exclude parts in the report
# Adds new "report_*" / optimize configuration rows to the "configurations" sheet,
# resizes its columns to fit the new content, updates sheet selections, and makes
# "configurations" the active tab again.

$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("configurations")
$wsDMO    = $wb.Worksheets.Item("decision_makers_options")

function Set-TextValue($ws, $addr, $text) {
    # Plain Range.Value coerces the literal strings "True"/"False" into real
    # booleans (t="b"). The source file stores them as ordinary shared-string
    # text, so those two values are entered as a formula and then flattened
    # back to a static value via copy / paste-special (values only), which
    # keeps them as plain text (t="s") without picking up a new cell style.
    $r = $ws.Range($addr)
    if ($text -eq "True" -or $text -eq "False") {
        $r.Formula = '="' + $text + '"'
        $r.Copy()
        $r.PasteSpecial(-4163)
    } else {
        $r.Value = $text
    }
}

# New configuration rows (3-12) on the "configurations" sheet.
Set-TextValue $wsConfig "A3"  "Optimize_DMO_name"
Set-TextValue $wsConfig "B3"  "Optimized_DMO"

Set-TextValue $wsConfig "A4"  "report_title_page"
Set-TextValue $wsConfig "B4"  "True"

Set-TextValue $wsConfig "A5"  "report_strategic_challenge"
Set-TextValue $wsConfig "B5"  "True"

Set-TextValue $wsConfig "A6"  "report_key_outputs_theme"
Set-TextValue $wsConfig "B6"  "True"

Set-TextValue $wsConfig "A7"  "report_decision_makers_options"
Set-TextValue $wsConfig "B7"  "True"

Set-TextValue $wsConfig "A8"  "report_scenarios"
Set-TextValue $wsConfig "B8"  "True"

Set-TextValue $wsConfig "A9"  "report_fixed_inputs"
Set-TextValue $wsConfig "B9"  "True"

Set-TextValue $wsConfig "A10" "report_dependencies"
Set-TextValue $wsConfig "B10" "False"

Set-TextValue $wsConfig "A11" "report_weighted_appreciations"
Set-TextValue $wsConfig "B11" "True"

Set-TextValue $wsConfig "A12" "report_add_optimize"
Set-TextValue $wsConfig "B12" "False"

# Widen the columns so the new, longer configuration names fit (mirrors the
# author resizing A:B after typing the new rows).
$wsConfig.Columns("A").ColumnWidth = 27.5
$wsConfig.Columns("B").ColumnWidth = 14

# Update the selection left on "decision_makers_options" (not the active
# sheet any more) before switching back to "configurations".
$wsDMO.Range("E5").Select()

# "configurations" becomes the active sheet/tab again, with its selection
# parked on E6.
$wsConfig.Activate()
$wsConfig.Range("E6").Select()
